# Insert a new row for an additional responsible teacher right after the
# existing "Docentes responsáveis:" entry (row 13), pushing every
# subsequent row down by one. Only columns B and C get the new value —
# column A keeps its original (empty) state for this row, matching the
# layout used by the other "Docentes responsáveis:" line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 14 (and everything below it) down, inserting a blank row 14.
$ws.Rows.Item(14).Insert()

# Fill in the new teacher's name in columns B and C of the freshly
# inserted row 14.
$ws.Cells.Item(14, 2).Value() = "1341653 - Maria José Ramos Sandim"
$ws.Cells.Item(14, 3).Value() = "1341653 - Maria José Ramos Sandim"
